# Update cryptocurrency price/volume data (and shifted Coin/Link entries for rows 11-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'28.197.31"
$ws.Range("E2").Value = '  +0.87%  '

# Row 3
$ws.Range("D3").Value = "'1.878.81"
$ws.Range("E3").Value = '  +1.23%  '

# Row 4
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").Value = "'313.78"
$ws.Range("E5").Value = '  +0.70%  '

# Row 6
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = '  +0.09%  '

# Row 7
$ws.Range("D7").Value = "'0.5125"
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").Value = "'0.3909"
$ws.Range("E8").Value = '  +3.07%  '

# Row 9
$ws.Range("D9").Value = "'0.08339"
$ws.Range("E9").Value = '  +0.73%  '

# Row 10
$ws.Range("D10").Value = "'1.120"
$ws.Range("E10").Value = '  +1.40%  '

# Row 11
$ws.Range("B11").Value = 'Polkadot'
$ws.Range("C11").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D11").Value = "'6.213"
$ws.Range("E11").Value = '  +0.53%  '

# Row 12
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = "'20.68"
$ws.Range("E12").Value = '  +1.38%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = "'1.882.80"
$ws.Range("E13").Value = '  +1.35%  '

# Row 14
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = "'7.279"
$ws.Range("E14").Value = '  +1.54%  '

# Row 15
$ws.Range("B15").Value = 'BinanceUSD'
$ws.Range("C15").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D15").Value = "'1.002"
$ws.Range("E15").Value = '  -0.19%  '

# Row 16
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = "'0.00001103"
$ws.Range("E16").Value = '  +0.91%  '

# Row 17
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = "'91.25"
$ws.Range("E17").Value = '  +1.18%  '

# Row 18
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = "'0.06652"
$ws.Range("E18").Value = '  +0.47%  '

# Row 19
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = "'17.78"
$ws.Range("E19").Value = '  +0.66%  '

# Row 20
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = '  +0.09%  '

# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = "'6.035"
$ws.Range("E21").Value = '  +0.58%  '

# Row 22
$ws.Range("B22").Value = 'WrappedBTC'
$ws.Range("C22").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D22").Value = "'28.240.04"
$ws.Range("E22").Value = '  +0.88%  '

# Row 23
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = "'11.16"
$ws.Range("E23").Value = '  +1.42%  '

# Row 24
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = "'2.267"
$ws.Range("E24").Value = '  +0.45%  '

# Row 25
$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").Value = "'2.095.45"
$ws.Range("E25").Value = '  +1.09%  '

# Row 26
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = "'2.512"
$ws.Range("E26").Value = '  -1.89%  '

# Row 27
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = "'158.92"
$ws.Range("E27").Value = '  +1.05%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = "'20.62"
$ws.Range("E28").Value = '  +1.29%  '

# Row 29
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = "'125.27"
$ws.Range("E29").Value = '  +0.68%  '

# Row 30
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = "'0.1065"
$ws.Range("E30").Value = '  +0.84%  '

# Row 31
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = "'1.043"
$ws.Range("E31").Value = '  +0.52%  '

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'5.874"
$ws.Range("E32").Value = '  +5.38%  '

# Row 33
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = "'3.596"
$ws.Range("E33").Value = '  -0.27%  '

# Row 34
$ws.Range("B34").Value = 'FraxShare'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D34").Value = "'9.719"
$ws.Range("E34").Value = '  +0.82%  '

# Row 35
$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").Value = "'0.02461"
$ws.Range("E35").Value = '  +2.29%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = "'0.06553"
$ws.Range("E36").Value = '  +0.67%  '

# Row 37
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").Value = "'0.2188"
$ws.Range("E37").Value = '  +1.76%  '

# Row 38
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").Value = "'1.206"
$ws.Range("E38").Value = '  +0.17%  '

# Row 39
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = "'0.6509"
$ws.Range("E39").Value = '  +1.98%  '

# Row 40
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").Value = "'5.018"
$ws.Range("E40").Value = '  +3.29%  '

# Row 41
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = "'1.230"
$ws.Range("E41").Value = '  +0.47%  '

# Row 42
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = "'11.30"
$ws.Range("E42").Value = '  +0.76%  '

# Row 43
$ws.Range("B43").Value = 'Decentraland'
$ws.Range("C43").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D43").Value = "'0.6133"
$ws.Range("E43").Value = '  +1.64%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = "'13.15"
$ws.Range("E44").Value = '  +0.64%  '

# Row 45
$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = "'1.286"
$ws.Range("E45").Value = '  +0.15%  '

# Row 46
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = "'3.681"
$ws.Range("E46").Value = '  +0.66%  '

# Row 47
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = "'2.014"
$ws.Range("E47").Value = '  +2.40%  '

# Row 48
$ws.Range("B48").Value = 'EOS'
$ws.Range("C48").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D48").Value = "'1.234"
$ws.Range("E48").Value = '  +2.48%  '

# Row 49
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = "'121.40"
$ws.Range("E49").Value = '  +0.92%  '

# Row 50
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = "'78.94"
$ws.Range("E50").Value = '  -0.68%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'0.06886"
$ws.Range("E51").Value = '  +0.83%  '
